# Handles float input without breaking stuff
#
# The marksheet previously rendered as "Absent" (a template/placeholder with
# 56 slots across 3 side-by-side question blocks, all zeroed). The real
# grading data for this student is now filled in: 28 actual questions,
# scored 17 right / 5 wrong / 6 not-attempted, with the redundant third
# question-block (columns G:H) removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary rows (9-12) -----------------------------------------------
# Give the row-label cells in column A the same title style used by their
# row-9 header neighbours.
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"

# Right / Wrong / Not-Attempt / Max
$ws.Range("B10").Value = 17
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 28

# Marking scheme: +4 for correct, -1 for incorrect (now numeric, not text)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Totals: 17*4 = 68, 5*-1 = -5, score 63 out of 112 (28 questions * 4 marks)
$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -5
$ws.Range("E12").Value = "63/112"

# --- Drop the third (G:H) question-answer block -------------------------
# Only two blocks (A:B and D:E) are still in use; removing the columns
# shrinks the used range from A5:H40 down to A5:E40 automatically.
$ws.Columns("G:H").Delete()

# --- Student answers (column A) for each question row -------------------
# Correct answers the student matched go in with "correctStyle" (green),
# mismatches get "incorrectStyle" (red), and un-attempted questions are
# left blank with "normalStyle".
$studentAnswers = @{
    16 = @{ Value = "Option A"; Style = "correctStyle" }
    17 = @{ Value = "Option D"; Style = "correctStyle" }
    18 = @{ Value = "Option B"; Style = "correctStyle" }
    19 = @{ Value = "Option C"; Style = "correctStyle" }
    20 = @{ Value = "Option B"; Style = "correctStyle" }
    21 = @{ Value = "Option C"; Style = "correctStyle" }
    23 = @{ Value = "Option D"; Style = "correctStyle" }
    25 = @{ Value = "Option A"; Style = "correctStyle" }
    26 = @{ Value = "Option C"; Style = "correctStyle" }
    27 = @{ Value = "Option C"; Style = "incorrectStyle" }
    28 = @{ Value = "Option C"; Style = "incorrectStyle" }
    31 = @{ Value = "Option D"; Style = "correctStyle" }
    33 = @{ Value = "Option D"; Style = "correctStyle" }
    35 = @{ Value = "Option D"; Style = "correctStyle" }
    36 = @{ Value = "Option A"; Style = "correctStyle" }
    37 = @{ Value = "Option A"; Style = "correctStyle" }
    38 = @{ Value = "Option B"; Style = "incorrectStyle" }
    39 = @{ Value = "Option C"; Style = "incorrectStyle" }
    40 = @{ Value = "Option D"; Style = "correctStyle" }
}

foreach ($row in $studentAnswers.Keys) {
    $entry = $studentAnswers[$row]
    $cell = $ws.Range("A$row")
    $cell.Value = $entry.Value
    $cell.Style = $entry.Style
}
# Rows 22, 24, 29, 30, 32, 34 stay blank/un-attempted (already normalStyle).

# --- Second question-answer block (D:E), only still used on rows 16-18 --
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"

$ws.Range("D17").Value = "Option A"
$ws.Range("D17").Style = "incorrectStyle"

$ws.Range("D18").Value = "Option D"
$ws.Range("D18").Style = "correctStyle"

# Every other row's D:E pair (19-40) is no longer used, so clear it back to
# blank while keeping formatting.
$ws.Range("D19:E40").ClearContents()
